$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-use the existing date-number-format style (already applied to D2:D10) for D1,
# then retype D1's header text from "Date " to "date".
$ws.Cells.Item(2, 4).Copy()
$ws.Cells.Item(1, 4).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(1, 4).Value = "date"

# Update the D column data values (date serial changes from 43624 to 36714, i.e. 2000-07-07)
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 4).Value = 36714
}

# Move the active selection to E1
$ws.Range("E1").Select()
